# Auto-generated Excel COM-interop script to apply the diff changes
# to the Sagittarius_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2221.3125
$ws.Range("I15").Value = 2221.3125
$ws.Range("K15").Value = 6663.9375
$ws.Range("M15").Value = -6494.9375

$ws.Range("H40").Value = 1681.0312
$ws.Range("I40").Value = 1719.72
$ws.Range("J40").Value = 1542.8572
$ws.Range("K40").Value = 1719.72
$ws.Range("L40").Value = 1542.8572
$ws.Range("M40").Value = -1544.72
$ws.Range("N40").Value = -1892.8572

$ws.Range("H100").Value = 1680.2
$ws.Range("I100").Value = 950
$ws.Range("J100").Value = 2167
$ws.Range("K100").Value = 950
$ws.Range("L100").Value = 2167
$ws.Range("M100").Value = -409
$ws.Range("N100").Value = -3249

$ws.Range("H107").Value = 2000.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2000.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2000.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5840.5

$ws.Range("H112").Value = 3322.3684
$ws.Range("I112").Value = 2794.5
$ws.Range("K112").Value = 8383.5
$ws.Range("M112").Value = -7275.5

$ws.Range("H132").Value = 2638.1667
$ws.Range("I132").Value = 2505.8
$ws.Range("K132").Value = 7517.400000000001
$ws.Range("M132").Value = -4987.400000000001

$ws.Range("H138").Value = 3966.5435
$ws.Range("I138").Value = 3067.182
$ws.Range("J138").Value = 4249.2
$ws.Range("K138").Value = 9201.545999999998
$ws.Range("L138").Value = 12747.6
$ws.Range("M138").Value = -4061.545999999998
$ws.Range("N138").Value = -23027.6

$ws.Range("H141").Value = 3267
$ws.Range("I141").Value = 2165
$ws.Range("J141").Value = 4920
$ws.Range("K141").Value = 6495
$ws.Range("L141").Value = 14760
$ws.Range("M141").Value = -1315
$ws.Range("N141").Value = -25120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 988.6
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 988.6
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H110").Value = 2282.3333
$ws.Range("I110").Value = 2282.3333
$ws.Range("K110").Value = 2282.3333
$ws.Range("M110").Value = -237.3332999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 15
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 163
$ws.Range("N22").Value = -366

$ws.Range("H80").Value = 181.4
$ws.Range("I80").Value = 174.66667
$ws.Range("J80").Value = 191.5
$ws.Range("K80").Value = 174.66667
$ws.Range("L80").Value = 191.5
$ws.Range("M80").Value = 823.3333299999999
$ws.Range("N80").Value = -2187.5

$ws.Range("H83").Value = 181.4
$ws.Range("I83").Value = 174.66667
$ws.Range("J83").Value = 191.5
$ws.Range("K83").Value = 873.3333500000001
$ws.Range("L83").Value = 957.5
$ws.Range("M83").Value = 4118.66665
$ws.Range("N83").Value = -10941.5

$ws.Range("H94").Value = 1894.0667
$ws.Range("I94").Value = 1851.1428
$ws.Range("K94").Value = 1851.1428
$ws.Range("M94").Value = -1400.1428

$ws.Range("H99").Value = 2480.3635
$ws.Range("I99").Value = 2480.3635
$ws.Range("K99").Value = 2480.3635
$ws.Range("M99").Value = -982.3634999999999

$ws.Range("H107").Value = 7400
$ws.Range("I107").Value = 7400
$ws.Range("K107").Value = 7400
$ws.Range("M107").Value = -5480

$ws.Range("H134").Value = 6147.25
$ws.Range("I134").Value = 6297
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 18891
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -16356
$ws.Range("N134").Value = -18570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7604.2
$ws.Range("I16").Value = 6007
$ws.Range("K16").Value = 6007
$ws.Range("M16").Value = -5720

$ws.Range("H113").Value = 7604.2
$ws.Range("I113").Value = 6007
$ws.Range("K113").Value = 6007
$ws.Range("M113").Value = -3837

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4314.4707
$ws.Range("J68").Value = 4314.4707
$ws.Range("L68").Value = 12943.4121
$ws.Range("N68").Value = -14565.4121

$ws.Range("H71").Value = 4314.4707
$ws.Range("J71").Value = 4314.4707
$ws.Range("L71").Value = 38830.2363
$ws.Range("N71").Value = -46942.2363

$ws.Range("H76").Value = 10500
$ws.Range("I76").Value = 10000
$ws.Range("J76").Value = 11000
$ws.Range("K76").Value = 30000
$ws.Range("L76").Value = 33000
$ws.Range("M76").Value = -29617
$ws.Range("N76").Value = -33766

$ws.Range("H79").Value = 10500
$ws.Range("I79").Value = 10000
$ws.Range("J79").Value = 11000
$ws.Range("K79").Value = 30000
$ws.Range("L79").Value = 33000
$ws.Range("M79").Value = -28674
$ws.Range("N79").Value = -35652

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

$ws.Range("H114").Value = 75.666664
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H131").Value = 478646.44
$ws.Range("J131").Value = 478646.44
$ws.Range("L131").Value = 1435939.32
$ws.Range("N131").Value = -1446019.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 78624.25
$ws.Range("J57").Value = 78624.25
$ws.Range("L57").Value = 78624.25
$ws.Range("N57").Value = -80264.25

$ws.Range("H102").Value = 3512.75
$ws.Range("I102").Value = 2748
$ws.Range("J102").Value = 4277.5
$ws.Range("K102").Value = 2748
$ws.Range("L102").Value = 4277.5
$ws.Range("M102").Value = -1126
$ws.Range("N102").Value = -7521.5

$ws.Range("H126").Value = 4858.923
$ws.Range("I126").Value = 3451.5
$ws.Range("J126").Value = 6065.2856
$ws.Range("K126").Value = 10354.5
$ws.Range("L126").Value = 18195.8568
$ws.Range("M126").Value = -7884.5
$ws.Range("N126").Value = -23135.8568

$ws.Range("H132").Value = 3926.2222
$ws.Range("I132").Value = 3926.2222
$ws.Range("K132").Value = 11778.6666
$ws.Range("M132").Value = -9248.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3574.5
$ws.Range("I22").Value = 2950
$ws.Range("J22").Value = 4199
$ws.Range("K22").Value = 2950
$ws.Range("L22").Value = 4199
$ws.Range("M22").Value = -2655
$ws.Range("N22").Value = -4789

$ws.Range("H27").Value = 3574.5
$ws.Range("I27").Value = 2950
$ws.Range("J27").Value = 4199
$ws.Range("K27").Value = 2950
$ws.Range("L27").Value = 4199
$ws.Range("M27").Value = -2843
$ws.Range("N27").Value = -4413

$ws.Range("H40").Value = 2323.5557
$ws.Range("I40").Value = 1925
$ws.Range("J40").Value = 2821.75
$ws.Range("K40").Value = 1925
$ws.Range("L40").Value = 2821.75
$ws.Range("M40").Value = -1789
$ws.Range("N40").Value = -3093.75

$ws.Range("H46").Value = 1671.5714
$ws.Range("J46").Value = 1399.5
$ws.Range("L46").Value = 1399.5
$ws.Range("N46").Value = -1775.5

$ws.Range("H61").Value = 2359.7144
$ws.Range("I61").Value = 2393
$ws.Range("J61").Value = 2315.3333
$ws.Range("K61").Value = 2393
$ws.Range("L61").Value = 2315.3333
$ws.Range("M61").Value = -2191
$ws.Range("N61").Value = -2719.3333

$ws.Range("H113").Value = 2359.7144
$ws.Range("I113").Value = 2393
$ws.Range("J113").Value = 2315.3333
$ws.Range("K113").Value = 2393
$ws.Range("L113").Value = 2315.3333
$ws.Range("M113").Value = -223
$ws.Range("N113").Value = -6655.3333

$ws.Range("H132").Value = 2251.5
$ws.Range("I132").Value = 2147.625
$ws.Range("J132").Value = 2667
$ws.Range("K132").Value = 6442.875
$ws.Range("L132").Value = 8001
$ws.Range("M132").Value = -3912.875
$ws.Range("N132").Value = -13061

$ws.Range("H136").Value = 2516.2
$ws.Range("I136").Value = 2350.75
$ws.Range("K136").Value = 7052.25
$ws.Range("M136").Value = -4502.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3273.3333
$ws.Range("I96").Value = 230
$ws.Range("J96").Value = 4795
$ws.Range("K96").Value = 230
$ws.Range("L96").Value = 4795
$ws.Range("M96").Value = 1143
$ws.Range("N96").Value = -7541

$ws.Range("H132").Value = 1900.2222
$ws.Range("I132").Value = 2202.1667
$ws.Range("K132").Value = 6606.500100000001
$ws.Range("M132").Value = -4076.500100000001
